$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 1531.1714
$ws.Range("I132").Value = 1341.4839
$ws.Range("J132").Value = 3001.25
$ws.Range("K132").Value = 4024.4517
$ws.Range("L132").Value = 9003.75
$ws.Range("M132").Value = -1494.4517
$ws.Range("N132").Value = -14063.75

# Row 138
$ws.Range("H138").Value = 2943.5476
$ws.Range("I138").Value = 1762.2667
$ws.Range("J138").Value = 3599.8147
$ws.Range("K138").Value = 5286.800099999999
$ws.Range("L138").Value = 10799.4441
$ws.Range("M138").Value = -146.8000999999995
$ws.Range("N138").Value = -21079.4441

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4787.25
$ws.Range("I61").Value = 3564.8
$ws.Range("K61").Value = 3564.8
$ws.Range("M61").Value = -3352.8

# Row 97
$ws.Range("H97").Value = 789.4545000000001
$ws.Range("I97").Value = 466
$ws.Range("J97").Value = 974.2857
$ws.Range("K97").Value = 466
$ws.Range("L97").Value = 974.2857
$ws.Range("M97").Value = 30
$ws.Range("N97").Value = -1966.2857

# Row 136
$ws.Range("H136").Value = 4787.25
$ws.Range("I136").Value = 3564.8
$ws.Range("K136").Value = 10694.4
$ws.Range("M136").Value = -8144.400000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 1831.25
$ws.Range("I64").Value = 1467
$ws.Range("J64").Value = 2091.4285
$ws.Range("K64").Value = 1467
$ws.Range("L64").Value = 2091.4285
$ws.Range("M64").Value = -1242
$ws.Range("N64").Value = -2541.4285

# Row 67
$ws.Range("H67").Value = 1831.25
$ws.Range("I67").Value = 1467
$ws.Range("J67").Value = 2091.4285
$ws.Range("K67").Value = 1467
$ws.Range("L67").Value = 2091.4285
$ws.Range("M67").Value = -687
$ws.Range("N67").Value = -3651.4285

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 1645
$ws.Range("I2").Value = 1645
$ws.Range("K2").Value = 1645
$ws.Range("M2").Value = -1532

# Row 22
$ws.Range("H22").Value = 4479.6
$ws.Range("I22").Value = 3532.6667
$ws.Range("J22").Value = 5900
$ws.Range("K22").Value = 3532.6667
$ws.Range("L22").Value = 5900
$ws.Range("M22").Value = -3182.6667
$ws.Range("N22").Value = -6600

# Row 31
$ws.Range("H31").Value = 56033.668
$ws.Range("I31").Value = 4467.3
$ws.Range("J31").Value = 102912.18
$ws.Range("K31").Value = 4467.3
$ws.Range("L31").Value = 102912.18
$ws.Range("M31").Value = -4172.3
$ws.Range("N31").Value = -103502.18

# Row 34
$ws.Range("H34").Value = 56033.668
$ws.Range("I34").Value = 4467.3
$ws.Range("J34").Value = 102912.18
$ws.Range("K34").Value = 4467.3
$ws.Range("L34").Value = 102912.18
$ws.Range("M34").Value = -4265.3
$ws.Range("N34").Value = -103316.18

# Row 135
$ws.Range("H135").Value = 68994.39999999999
$ws.Range("J135").Value = 68994.39999999999
$ws.Range("L135").Value = 68994.39999999999
$ws.Range("N135").Value = -79134.39999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 887.3333
$ws.Range("I14").Value = 887.3333
$ws.Range("K14").Value = 2661.9999
$ws.Range("M14").Value = -2488.9999

# Row 40
$ws.Range("H40").Value = 352.75
$ws.Range("I40").Value = 21
$ws.Range("J40").Value = 400.14285
$ws.Range("K40").Value = 84
$ws.Range("L40").Value = 1600.5714
$ws.Range("M40").Value = -15
$ws.Range("N40").Value = -1738.5714

# Row 122
$ws.Range("H122").Value = 1680.72
$ws.Range("J122").Value = 1757.6086
$ws.Range("L122").Value = 15818.4774
$ws.Range("N122").Value = -20718.4774

# Row 137
$ws.Range("H137").Value = 5484
$ws.Range("I137").Value = 1743.4
$ws.Range("J137").Value = 9224.6
$ws.Range("K137").Value = 5230.200000000001
$ws.Range("L137").Value = 27673.8
$ws.Range("M137").Value = -130.2000000000007
$ws.Range("N137").Value = -37873.8

# Row 138
$ws.Range("H138").Value = 8652
$ws.Range("I138").Value = 5612.8
$ws.Range("K138").Value = 16838.4
$ws.Range("M138").Value = -11698.4

# Row 139
$ws.Range("H139").Value = 3767.8696
$ws.Range("I139").Value = 1668.2778
$ws.Range("J139").Value = 11326.4
$ws.Range("K139").Value = 5004.8334
$ws.Range("L139").Value = 33979.2
$ws.Range("M139").Value = 135.1665999999996
$ws.Range("N139").Value = -44259.2

$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 32944.273
$ws.Range("J93").Value = 33415.9
$ws.Range("L93").Value = 33415.9
$ws.Range("N93").Value = -37159.9

# Row 132
$ws.Range("H132").Value = 73998.625
$ws.Range("I132").Value = 134123
$ws.Range("K132").Value = 402369
$ws.Range("M132").Value = -399839

$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

# Row 22
$ws.Range("H22").Value = 5800.4
$ws.Range("I22").Value = 2500
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = -2205

# Row 27
$ws.Range("H27").Value = 5800.4
$ws.Range("I27").Value = 2500
$ws.Range("K27").Value = 2500
$ws.Range("M27").Value = -2393

# Row 28
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

# Row 37
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

# Row 61
$ws.Range("H61").Value = 11626.625
$ws.Range("I61").Value = 10858.286
$ws.Range("J61").Value = 17005
$ws.Range("K61").Value = 10858.286
$ws.Range("L61").Value = 17005
$ws.Range("M61").Value = -10656.286
$ws.Range("N61").Value = -17409

# Row 93
$ws.Range("H93").Value = 2536.7273
$ws.Range("I93").Value = 1790.9
$ws.Range("K93").Value = 1790.9
$ws.Range("M93").Value = -542.9000000000001

# Row 113
$ws.Range("H113").Value = 11626.625
$ws.Range("I113").Value = 10858.286
$ws.Range("J113").Value = 17005
$ws.Range("K113").Value = 10858.286
$ws.Range("L113").Value = 17005
$ws.Range("M113").Value = -8688.286
$ws.Range("N113").Value = -21345

$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 4000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 4000
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -4228

# Row 5
$ws.Range("H5").Value = 1502001
$ws.Range("J5").Value = 1502001
$ws.Range("L5").Value = 1502001
$ws.Range("N5").Value = -1502225

Write-Output "edits applied"